$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4 with corrected per-mandataire taxe/brut amounts ---
$ws.Range("G2").Value = 4000
$ws.Range("J2").Value = 400
$ws.Range("M2").Value = 3600

$ws.Range("G3").Value = 3000
$ws.Range("J3").Value = 300
$ws.Range("M3").Value = 2700

$ws.Range("G4").Value = 3000
$ws.Range("J4").Value = 300
$ws.Range("M4").Value = 2700

# --- Add new rows 5-8 ---
$ws.Range("A5").Value = "001/LF/DR IFRAN"
$ws.Range("B5").Value = "Logement de fonction"
$ws.Range("C5").Value = "BB12354"
$ws.Range("D5").Value = "Tawfiq MF"
$ws.Range("E5").Value = "ds"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 1000
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1000
$ws.Range("M5").Value = 1000

$ws.Range("A6").Value = "001/LF/DR IFRAN"
$ws.Range("B6").Value = "Logement de fonction"
$ws.Range("C6").Value = "BB123456"
$ws.Range("D6").Value = "Ahmed tawfiq"
$ws.Range("E6").Value = "ds"
$ws.Range("F6").Value = "mensuelle"
$ws.Range("G6").Value = 2000
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000

$ws.Range("A7").Value = "001/LF/DR IFRAN"
$ws.Range("B7").Value = "Logement de fonction"
$ws.Range("C7").Value = "bs3"
$ws.Range("D7").Value = "IBM"
$ws.Range("E7").Value = "ds"
$ws.Range("F7").Value = "mensuelle"
$ws.Range("G7").Value = 3000
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 10
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 2700

$ws.Range("A8").Value = "001/LF/DR IFRAN"
$ws.Range("B8").Value = "Logement de fonction"
$ws.Range("C8").Value = "BB123459"
$ws.Range("D8").Value = "mamadu sacko"
$ws.Range("E8").Value = "ds"
$ws.Range("F8").Value = "mensuelle"
$ws.Range("G8").Value = 4000
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 10
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = 3600
